$d = $word.ActiveDocument

# 1) "QUE EL PROMITENTE VENDEDOR LE ENTREGUE" -> "QUE "{{SEXO_1}} PROMITENTE {{SEXO_2}}" LE ENTREGUE"
$d.Content.Find.Execute(
    "EL PROMITENTE VENDEDOR",
    $true, $false, $false, $false, $false, $true, 1, $false,
    [char]0x201C + "{{SEXO_1}} PROMITENTE {{SEXO_2}}" + [char]0x201D,
    2)

# 2) Pluralize verbs: INCURRA -> INCURRAN, COMUNIQUE -> COMUNIQUEN, SOLICITE -> SOLICITEN
$d.Content.Find.Execute(
    "INCURRA EN EL INCUMPLIMIENTO",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "INCURRAN EN EL INCUMPLIMIENTO",
    2)

$d.Content.Find.Execute(
    "NO SE COMUNIQUE NI SOLICITE POR ESCRITO",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "NO SE COMUNIQUEN NI SOLICITEN POR ESCRITO",
    2)

Write-Output "done"
